# Updates the cryptos price/volume snapshot (and two coin-rank swaps)
# to match the latest scrape, per the commit
# "Updated cryptos list on Sat Mar  9 22:32:53 UTC 2024 with GitHub Actions".
#
# Cells in column D/E are plain text (prices use dotted thousand
# separators like '68.365.77', percentages are padded strings like
# '  +0.20%  '), so numeric-looking values are written with a leading
# apostrophe to keep Excel from re-typing them as numbers (which would
# drop meaningful trailing zeros, e.g. 487.20 -> 487.2). The cell style
# is reset to Normal right after so the quote-prefix flag introduced by
# the apostrophe doesn't leave a spurious style change behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.365.77"
$ws.Range("E2").Value = "  +0.20%  "
# Row 3
$ws.Range("D3").Value = "3.904.70"
$ws.Range("E3").Value = "  +0.01%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").Value = "'487.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
# Row 6
$ws.Range("D6").Value = "'146.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "
# Row 7
$ws.Range("E7").Value = "  +0.16%  "
# Row 8
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
$ws.Range("D9").Value = "'0.745"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.99%  "
# Row 10
$ws.Range("E10").Value = "  +9.31%  "
# Row 11
$ws.Range("D11").Value = "'0.0000358"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
# Row 12
$ws.Range("D12").Value = "'43.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.56%  "
# Row 13
$ws.Range("D13").Value = "'10.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.79%  "
# Row 14
$ws.Range("D14").Value = "4.520.92"
$ws.Range("E14").Value = "  -0.17%  "
# Row 15
$ws.Range("D15").Value = "3.932.54"
$ws.Range("E15").Value = "  +0.99%  "
# Row 16
$ws.Range("E16").Value = "  -2.37%  "
# Row 17
$ws.Range("E17").Value = "  -0.51%  "
# Row 18
$ws.Range("D18").Value = "'20.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "
# Row 19
$ws.Range("E19").Value = "  +1.99%  "
# Row 20
$ws.Range("D20").Value = "68.439.01"
$ws.Range("E20").Value = "  +0.15%  "
# Row 21
$ws.Range("D21").Value = "'432.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "
# Row 22
$ws.Range("E22").Value = "  +5.58%  "
# Row 23
$ws.Range("E23").Value = "  +0.75%  "
# Row 24
$ws.Range("D24").Value = "'89.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "
# Row 25
$ws.Range("D25").Value = "'12.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +17.73%  "
# Row 26
$ws.Range("D26").Value = "'3.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.74%  "
# Row 27
$ws.Range("D27").Value = "'10.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.11%  "
# Row 28
$ws.Range("D28").Value = "'37.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.32%  "
# Row 29
$ws.Range("D29").Value = "'5.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "
# Row 30
$ws.Range("D30").Value = "'712.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.35%  "
# Row 31
$ws.Range("D31").Value = "'0.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "
# Row 32
$ws.Range("D32").Value = "'13.41"
$ws.Range("D32").Style = "Normal"
# Row 33
$ws.Range("E33").Value = "  +2.26%  "
# Row 34
$ws.Range("D34").Value = "0.0₃0901"
$ws.Range("E34").Value = "  -2.09%  "
# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.10%  "
# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'61.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.76%  "
# Row 37
$ws.Range("D37").Value = "'40.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "
# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.404"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.87%  "
# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.77%  "
# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
# Row 41
$ws.Range("E41").Value = "  +4.72%  "
# Row 42
$ws.Range("D42").Value = "'2.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.79%  "
# Row 43
$ws.Range("E43").Value = "  +3.02%  "
# Row 44
$ws.Range("E44").Value = "  -0.46%  "
# Row 45
$ws.Range("D45").Value = "0.0₆0376"
$ws.Range("E45").Value = "  +28.76%  "
# Row 46
$ws.Range("E46").Value = "  +1.04%  "
# Row 47
$ws.Range("D47").Value = "'3.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.83%  "
# Row 48
$ws.Range("E48").Value = "  +0.03%  "
# Row 49
$ws.Range("E49").Value = "  -1.57%  "
# Row 50
$ws.Range("D50").Value = "'2.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "
# Row 51
$ws.Range("D51").Value = "'142.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.18%  "
